$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 8): "PRÓXIMO SETOR" column removed, "META HT" column
#    added just before "HT". This shifts Nº O.S / QTD / STATUS / PENDÊNCIA
#    one column to the right-ish and extends the "PRODUTO SERVIÇO" merge.
# ---------------------------------------------------------------------------

# Break apart the old merged groups in row 8 before rebuilding them.
$ws.Range("G8:K8").UnMerge()
$ws.Range("N8:O8").UnMerge()
$ws.Range("P8:Q8").UnMerge()
$ws.Range("R8:S8").UnMerge()

# Shift the trailing header labels into their new homes.
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = "Nº O.S"
$ws.Range("N8").Value = "QTD"
$ws.Range("O8").Value = "STATUS"
$ws.Range("P8").Value = ""
$ws.Range("Q8").Value = "PENDÊNCIA E/OU OBSERVAÇÃO"
$ws.Range("R8").Value = ""
$ws.Range("S8").Value = "META HT"

# Re-merge "PRODUTO SERVIÇO" across one extra column (now ends at L instead of K).
$ws.Range("G8:L8").Merge()
$ws.Range("O8:P8").Merge()
$ws.Range("Q8:R8").Merge()

# ---------------------------------------------------------------------------
# Formatting for the rebuilt header cells: same look as the rest of row 8
# (8pt font, centered, thin box borders) - copy from the existing matching
# box-border styles already used elsewhere on the sheet, then touch up the
# font/alignment so it reads the same as before.
# ---------------------------------------------------------------------------

# Reset every touched cell back to the plain single-cell boxed-header look
# first (same style as B8:F8/T8), then special-case the merged groups below.
$ws.Range("T8").Copy()
$ws.Range("G8:S8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Left edge of a merged box (border: left+top+bottom)
$ws.Range("H4").Copy()
$leftCells = $ws.Range("G8,O8,Q8")
$leftCells.PasteSpecial(-4122)
$leftCells.Font.Size = 8
$leftCells.VerticalAlignment = -4108

# Middle of a merged box (border: top+bottom only)
$ws.Range("I4").Copy()
$midCells = $ws.Range("H8:K8")
$midCells.PasteSpecial(-4122)
$midCells.Font.Size = 8
$midCells.VerticalAlignment = -4108

# Right edge of a merged box (border: right+top+bottom)
$ws.Range("M3").Copy()
$rightCells = $ws.Range("L8,P8,R8")
$rightCells.PasteSpecial(-4122)
$rightCells.Font.Size = 8
$rightCells.HorizontalAlignment = -4108

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) New column width for the inserted "META HT" column (P).
# ---------------------------------------------------------------------------
$ws.Columns.Item(16).ColumnWidth = 9.17

# ---------------------------------------------------------------------------
# 3) Selection moves to O16 (cosmetic, matches the saved cursor position).
# ---------------------------------------------------------------------------
$ws.Range("O16").Select()
